$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Add "#230: fix bug: arrows animation" as a new run into the
#    currently-empty bullet paragraph right after
#    "#235 Check/Fix DeleteSIgGenCommand (not reproduceable)"
#    (the "fixed in this version" list near the top of the doc).
# ------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^#235 Check/Fix DeleteSIgGenCommand") {
        $targetPara = $p.Next()
        break
    }
}

$r = $targetPara.Range
$r.InsertAfter("#230: fix bug: arrows animation")
$r2 = $targetPara.Range
$r2.Font.Name = "Arial"
$r2.Font.NameFarEast = "Times New Roman"
$r2.Font.NameBi = "Arial"
$r2.Font.Size = 9
$r2.Font.SizeBi = 9
$r2.Font.Color = 3025188
$r2.LanguageID = "en-US"
$r2.LanguageIDFarEast = "de-DE"

# ------------------------------------------------------------------
# 2) Split "#162 Fix DEF_FUNC (new without delete)" so that the
#    function name becomes its own run reading "ScrDefConst"
#    (identical formatting, just forced onto a separate run).
# ------------------------------------------------------------------
$defFuncPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "DEF_FUNC") {
        $defFuncPara = $p
        break
    }
}

$full = $defFuncPara.Range.Text
$idx = $full.IndexOf("DEF_FUNC")
$subStart = $defFuncPara.Range.Start + $idx
$subEnd = $subStart + 8
$sub = $d.Range($subStart, $subEnd)
$sub.Text = "ScrDefConst"

# Force the run boundary to actually split (even though the
# resulting formatting is identical to its neighbours) by toggling a
# character attribute on and back off.
$subEnd2 = $subStart + 11
$sub2 = $d.Range($subStart, $subEnd2)
$sub2.Font.Bold = $true
$sub2.Font.Bold = $false

# ------------------------------------------------------------------
# 3) Remove the now-duplicated "#230: fix bug: arrows animation "
#    text (plus its trailing single-space run) from the bug table,
#    where it used to precede "#34 Restore last model version".
# ------------------------------------------------------------------
$tablePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^#230: fix bug: arrows animation #34 Restore last") {
        $tablePara = $p
        break
    }
}

$delStart = $tablePara.Range.Start
$delEnd = $delStart + 32
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()
